$wb = $excel.ActiveWorkbook

# Codebook sheet selection state change (done first so the Dictionary
# Mapping tab ends up as the active/selected tab, matching the author's
# final saved state)
$codebook = $wb.Worksheets.Item("Codebook")
$codebook.Range("A6").Select() | Out-Null

$ws = $wb.Worksheets.Item("Dictionary Mapping")

# Row 14: Relation (H14) sio:SIO_000068 -> sio:SIO_000059
$ws.Range("H14").Value = "sio:SIO_000059"

# Row 16: ??ucm now relates directly (H16 added) to urine, with sio:SIO_000059
$ws.Range("H16").Value = "sio:SIO_000059"

# Row 17: URXUCM - fix Attribute id + add Unit of Measure (D17)
$ws.Range("B17").Value = "sio:SIO_001088"
$ws.Range("D17").Value = "uo:0000301"

# Row 18: ??ins - add Relation (H18), fix inRelationTo (I18) blood instead of ??isn
$ws.Range("H18").Value = "sio:SIO_000059"
$ws.Range("I18").Value = "??blood"

# Row 19: Column renamed LBXIN -> LBDNISI, fix Attribute id, fix inRelationTo -> ??isn
$ws.Range("A19").Value = "LBDNISI"
$ws.Range("B19").Value = "sio:SIO_001088"
$ws.Range("I19").Value = "??isn"

# Row 20: ??uio - add Relation (H20), fix inRelationTo -> ??urine
$ws.Range("H20").Value = "sio:SIO_000059"
$ws.Range("I20").Value = "??urine"

# Row 21: WTSA2YR - fix Attribute id, add Unit of Measure (D21), fix inRelationTo -> ??uio
$ws.Range("B21").Value = "sio:SIO_001088"
$ws.Range("D21").Value = "uo:0000301"
$ws.Range("I21").Value = "??uio"

# Row 22: ??vid - add Relation (H22), fix inRelationTo -> ??blood
$ws.Range("H22").Value = "sio:SIO_000068"
$ws.Range("I22").Value = "??blood"

# Row 23: LBXVIDMS - fix Attribute id, add Unit of Measure (D23), fix inRelationTo -> ??vid
$ws.Range("B23").Value = "sio:SIO_001088"
$ws.Range("D23").Value = "uo:0000041"
$ws.Range("I23").Value = "??vid"

# Row 24 (??hepb) removed entirely
$ws.Range("A24:K24").ClearContents()

# Row 25: LBXHBC - fix Attribute id (was placeholder sio:SIO_), add inRelationTo
$ws.Range("B25").Value = "ncit:C75678"
$ws.Range("I25").Value = "??blood"

# Row 26: LBXBHS - add Attribute + inRelationTo
$ws.Range("B26").Value = "ncit:C628795"
$ws.Range("I26").Value = "??blood"

# Row 27 (new): LBDHBG
$ws.Range("A27").Value = "LBDHBG"
$ws.Range("B27").Value = "ncit:C75678"
$ws.Range("C27").Value = "??blood"
$ws.Range("I27").Value = "??blood"

# Row 28 (new): LBDHD
$ws.Range("A28").Value = "LBDHD"
$ws.Range("B28").Value = "ncit:C96664"
$ws.Range("C28").Value = "??blood"
$ws.Range("I28").Value = "??blood"

# Update selection state to match the author's final cursor position
# (Dictionary Mapping stays the active tab since this runs last)
$ws.Range("D27").Select() | Out-Null
